# WR_89877351_WeekEnding_071325.xlsx — single-WR enforcement edit
# Regenerated report: removes a line item from the Saturday section,
# refreshes totals/timestamp, and clears the Scope ID.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the "Point 14 / SVC-VA" line item row from the Saturday table.
#    Deleting the whole row shifts everything below it up by one and keeps
#    merged-cell references / the Sunday table in sync automatically.
$ws.Rows("17").Delete()

# 2) Refresh the "Report Generated On" timestamp.
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:02 AM"

# 3) Update the summary box: Total Billed Amount + Total Line Items.
$ws.Range("C8").Value = 732.61
$ws.Range("C9").Value = 3

# 4) Clear the Scope ID # value (now blank).
$ws.Range("G10").Value = ""

# 5) Saturday table — update the remaining line item price and its TOTAL.
$ws.Range("H16").Value = 198.88
$ws.Range("H17").Value = 198.88

# 6) Sunday table (rows shifted up by 1 after the delete above) — update
#    both line item prices and the TOTAL.
$ws.Range("H22").Value = 478.55
$ws.Range("H23").Value = 55.18
$ws.Range("H24").Value = 533.73
